$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for year 2025 (rows 16-22)
$data = @(
    @(2025, "Violencia doméstica", 1580),
    @(2025, "Acecho", 68),
    @(2025, "Agresión sexual", 24),
    @(2025, "Violencia en cita", 1),
    @(2025, "Discrimen de género", 78),
    @(2025, "Otras", 67),
    @(2025, "Trata Humana", 0)
)

$row = 16
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("B20").Select()
